# Supplyworks.xlsx - "Update info capture from master_url file"
#
# Changes applied:
#  - Fix accent typo in the "Homólogo Mansfield" header -> "Homologo Mansfield"
#  - Remove the hyperlink behind the Link cell (keep the URL text)
#  - Insert a new "Type" column (value "Tank") before "Linea"
#  - Rename the "Name" header to "Descripcion"
#  - Insert a new "Short Name" column (value "Gerber Ultra Flush 1 gpf Tank")
#    right after "Descripcion" and before "Link"
#  - Re-center the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the accented header text -------------------------------------
$ws.Range("B1").Value = "Homologo Mansfield"

# --- Drop the hyperlink on the Link cell but keep its visible text ----
$ws.Range("J2").Hyperlinks.Delete()

# --- Insert "Type" / "Tank" column before the "Linea" column (D) ------
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "Type"
$ws.Range("D2").Value = "Tank"
# match the width of its left neighbour (Sku) like the rest of the block
$ws.Columns("D:D").ColumnWidth = $ws.Columns("C:C").ColumnWidth

# --- "Name" header becomes "Descripcion" (now shifted to column J) ----
$ws.Range("J1").Value = "Descripcion"

# --- Insert "Short Name" column after "Descripcion", before "Link" ----
$ws.Columns("K:K").Insert()
$ws.Range("K1").Value = "Short Name"
$ws.Range("K2").Value = "Gerber Ultra Flush 1 gpf Tank"
$ws.Columns("K:K").ColumnWidth = 25.86

# --- Re-center every header cell (A1:L1) --------------------------------
$ws.Rows("1:1").HorizontalAlignment = -4108

# the new "Type" column starts a new visual group -> no left border,
# matching the border the "Linea" header used to have
$ws.Range("D1").Borders.Item(7).LineStyle = -4142

# --- Leave the selection where the author left it ----------------------
$ws.Range("J10").Select()
